$d = $word.ActiveDocument

$d.Content.Find.Execute('87×18=1566', $true, $false, $false, $false, $false, $true, 1, $false, '99×73=7227', 2) | Out-Null
$d.Content.Find.Execute('43×44=1892', $true, $false, $false, $false, $false, $true, 1, $false, '78×84=6552', 2) | Out-Null
$d.Content.Find.Execute('68×14=952', $true, $false, $false, $false, $false, $true, 1, $false, '22×95=2090', 2) | Out-Null
$d.Content.Find.Execute('15×13=195', $true, $false, $false, $false, $false, $true, 1, $false, '90×96=8640', 2) | Out-Null
$d.Content.Find.Execute('37×46=1702', $true, $false, $false, $false, $false, $true, 1, $false, '74×35=2590', 2) | Out-Null
$d.Content.Find.Execute('67×78=5226', $true, $false, $false, $false, $false, $true, 1, $false, '56×37=2072', 2) | Out-Null
$d.Content.Find.Execute('29×57=1653', $true, $false, $false, $false, $false, $true, 1, $false, '66×62=4092', 2) | Out-Null
$d.Content.Find.Execute('59×18=1062', $true, $false, $false, $false, $false, $true, 1, $false, '49×47=2303', 2) | Out-Null
$d.Content.Find.Execute('85×94=7990', $true, $false, $false, $false, $false, $true, 1, $false, '53×17=901', 2) | Out-Null
$d.Content.Find.Execute('69×67=4623', $true, $false, $false, $false, $false, $true, 1, $false, '56×96=5376', 2) | Out-Null
$d.Content.Find.Execute('41×13=533', $true, $false, $false, $false, $false, $true, 1, $false, '43×68=2924', 2) | Out-Null
$d.Content.Find.Execute('14×54=756', $true, $false, $false, $false, $false, $true, 1, $false, '66×76=5016', 2) | Out-Null
$d.Content.Find.Execute('22×76=1672', $true, $false, $false, $false, $false, $true, 1, $false, '14×22=308', 2) | Out-Null
$d.Content.Find.Execute('68×35=2380', $true, $false, $false, $false, $false, $true, 1, $false, '64×61=3904', 2) | Out-Null
$d.Content.Find.Execute('44×87=3828', $true, $false, $false, $false, $false, $true, 1, $false, '75×33=2475', 2) | Out-Null
$d.Content.Find.Execute('98×24=2352', $true, $false, $false, $false, $false, $true, 1, $false, '54×23=1242', 2) | Out-Null
$d.Content.Find.Execute('47×83=3901', $true, $false, $false, $false, $false, $true, 1, $false, '68×97=6596', 2) | Out-Null
$d.Content.Find.Execute('71×40=2840', $true, $false, $false, $false, $false, $true, 1, $false, '76×43=3268', 2) | Out-Null
$d.Content.Find.Execute('25×59=1475', $true, $false, $false, $false, $false, $true, 1, $false, '84×52=4368', 2) | Out-Null
$d.Content.Find.Execute('40×86=3440', $true, $false, $false, $false, $false, $true, 1, $false, '80×45=3600', 2) | Out-Null
$d.Content.Find.Execute('18×97=1746', $true, $false, $false, $false, $false, $true, 1, $false, '93×31=2883', 2) | Out-Null
$d.Content.Find.Execute('45×82=3690', $true, $false, $false, $false, $false, $true, 1, $false, '97×87=8439', 2) | Out-Null
$d.Content.Find.Execute('63×99=6237', $true, $false, $false, $false, $false, $true, 1, $false, '44×16=704', 2) | Out-Null
$d.Content.Find.Execute('79×15=1185', $true, $false, $false, $false, $false, $true, 1, $false, '36×38=1368', 2) | Out-Null
$d.Content.Find.Execute('72×17=1224', $true, $false, $false, $false, $false, $true, 1, $false, '91×32=2912', 2) | Out-Null
